{"js": "// Remove the whole \"\u6587\u4ef6\u8bf4\u660e\uff1a\" section (heading + its bulleted\n// description paragraphs + the trailing blank paragraph) that used to\n// sit at the very top of the document, right before the\n// \"\u8f6f\u4ef6\u5b89\u88c5\u8bf4\u660e\uff1a\" section. (commit: \"for PMT docs update\")\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the first paragraph of the section that must stay (\"\u8f6f\u4ef6\u5b89\u88c5\u8bf4\u660e\uff1a\").\n// Everything before it (the \"\u6587\u4ef6\u8bf4\u660e\uff1a\" heading, its numbered bullet\n// items, and the blank spacer paragraph) is removed.\nconst keepMarker = \"\u8f6f\u4ef6\u5b89\u88c5\u8bf4\u660e\uff1a\";\nlet keepIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === keepMarker) {\n    keepIndex = i;\n    break;\n  }\n}\n\nif (keepIndex === -1) {\n  throw new Error('Could not locate the \"' + keepMarker + '\" paragraph.');\n}\n\n// Delete paragraphs [0 .. keepIndex-1] in reverse order so indices stay\n// valid as items are removed.\nfor (let i = keepIndex - 1; i >= 0; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the whole \"\u6587\u4ef6\u8bf4\u660e\uff1a\" section (heading + its bulleted\n# description paragraphs + the trailing blank paragraph) that used to\n# sit at the very top of the document, right before the\n# \"\u8f6f\u4ef6\u5b89\u88c5\u8bf4\u660e\uff1a\" section. (commit: \"for PMT docs update\")\n\n$d = $word.ActiveDocument\n\n# Locate the start of the paragraph that must be kept (\"\u8f6f\u4ef6\u5b89\u88c5\u8bf4\u660e\uff1a\")\n# by searching the document content for its text.\n$finder = $d.Content\n$finder.Find.Execute(\"\u8f6f\u4ef6\u5b89\u88c5\u8bf4\u660e\") | Out-Null\n\nif (-not $finder.Find.Found) {\n    throw \"Could not locate the '\u8f6f\u4ef6\u5b89\u88c5\u8bf4\u660e\uff1a' paragraph.\"\n}\n\n# Build a range from the very start of the document up to the start of\n# that match (i.e. everything before it) and delete it in one shot.\n$deleteRange = $d.Range(0, $finder.Start)\n$deleteRange.Delete()\n"}
